$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values - repulled data / recalculated means
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 9
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = -4
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = 6
$ws.Range("F13").Value = -6
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = -1
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = -1
